# Auto-generated Excel COM-interop script
# Updates Leve profit calculation values (columns H-N) across all crafting job sheets
# per the scheduled market-data refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 127.542854
$ws.Range("I33").Value = 118.451614
$ws.Range("J33").Value = 198
$ws.Range("K33").Value = 118.451614
$ws.Range("L33").Value = 198
$ws.Range("M33").Value = 110.548386
$ws.Range("N33").Value = -656
$ws.Range("H70").Value = 13973373
$ws.Range("I70").Value = 33534194
$ws.Range("J70").Value = 1357.5714
$ws.Range("K70").Value = 100602582
$ws.Range("L70").Value = 4072.7142
$ws.Range("M70").Value = -100602312
$ws.Range("N70").Value = -4612.7142
$ws.Range("H73").Value = 13973373
$ws.Range("I73").Value = 33534194
$ws.Range("J73").Value = 1357.5714
$ws.Range("K73").Value = 100602582
$ws.Range("L73").Value = 4072.7142
$ws.Range("M73").Value = -100601646
$ws.Range("N73").Value = -5944.7142
$ws.Range("H98").Value = 3322.6155
$ws.Range("J98").Value = 7992.7144
$ws.Range("L98").Value = 7992.7144
$ws.Range("N98").Value = -10988.7144
$ws.Range("H122").Value = 3322.6155
$ws.Range("J122").Value = 7992.7144
$ws.Range("L122").Value = 23978.1432
$ws.Range("N122").Value = -28878.1432
$ws.Range("H137").Value = 1292.2142
$ws.Range("I137").Value = 1066.4667
$ws.Range("J137").Value = 1856.5834
$ws.Range("K137").Value = 3199.4001
$ws.Range("L137").Value = 5569.7502
$ws.Range("M137").Value = -649.4000999999998
$ws.Range("N137").Value = -10669.7502
$ws.Range("H138").Value = 4869.814
$ws.Range("I138").Value = 1601.9131
$ws.Range("K138").Value = 4805.7393
$ws.Range("M138").Value = 334.2606999999998
$ws.Range("H141").Value = 4546.558
$ws.Range("I141").Value = 2253.5854
$ws.Range("J141").Value = 51552.5
$ws.Range("K141").Value = 6760.7562
$ws.Range("L141").Value = 154657.5
$ws.Range("M141").Value = -1580.7562
$ws.Range("N141").Value = -165017.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8568.934999999999
$ws.Range("I32").Value = 9697.297
$ws.Range("J32").Value = 2551
$ws.Range("K32").Value = 9697.297
$ws.Range("L32").Value = 2551
$ws.Range("M32").Value = -9410.297
$ws.Range("N32").Value = -3125
$ws.Range("H61").Value = 995.55554
$ws.Range("I61").Value = 899.5946
$ws.Range("J61").Value = 1204.4117
$ws.Range("K61").Value = 899.5946
$ws.Range("L61").Value = 1204.4117
$ws.Range("M61").Value = -687.5946
$ws.Range("N61").Value = -1628.4117
$ws.Range("H63").Value = 4564.25
$ws.Range("I63").Value = 4524.636
$ws.Range("J63").Value = 5000
$ws.Range("K63").Value = 4524.636
$ws.Range("L63").Value = 5000
$ws.Range("M63").Value = -3838.636
$ws.Range("N63").Value = -6372
$ws.Range("H66").Value = 4564.25
$ws.Range("I66").Value = 4524.636
$ws.Range("J66").Value = 5000
$ws.Range("K66").Value = 22623.18
$ws.Range("L66").Value = 25000
$ws.Range("M66").Value = -19191.18
$ws.Range("N66").Value = -31864
$ws.Range("H74").Value = 676.5454999999999
$ws.Range("I74").Value = 649.2889
$ws.Range("J74").Value = 799.2
$ws.Range("K74").Value = 649.2889
$ws.Range("L74").Value = 799.2
$ws.Range("M74").Value = 224.7111
$ws.Range("N74").Value = -2547.2
$ws.Range("H77").Value = 676.5454999999999
$ws.Range("I77").Value = 649.2889
$ws.Range("J77").Value = 799.2
$ws.Range("K77").Value = 3246.4445
$ws.Range("L77").Value = 3996
$ws.Range("M77").Value = 1121.5555
$ws.Range("N77").Value = -12732
$ws.Range("H133").Value = 80729.75
$ws.Range("J133").Value = 80729.75
$ws.Range("L133").Value = 80729.75
$ws.Range("N133").Value = -85789.75
$ws.Range("H136").Value = 995.55554
$ws.Range("I136").Value = 899.5946
$ws.Range("J136").Value = 1204.4117
$ws.Range("K136").Value = 2698.7838
$ws.Range("L136").Value = 3613.2351
$ws.Range("M136").Value = -148.7838000000002
$ws.Range("N136").Value = -8713.2351

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 39462
$ws.Range("J35").Value = 39462
$ws.Range("L35").Value = 39462
$ws.Range("N35").Value = -40082
$ws.Range("H82").Value = 270000
$ws.Range("I82").Value = 270000
$ws.Range("K82").Value = 270000
$ws.Range("M82").Value = -269617
$ws.Range("H85").Value = 270000
$ws.Range("I85").Value = 270000
$ws.Range("K85").Value = 270000
$ws.Range("M85").Value = -268674
$ws.Range("H134").Value = 1525.585
$ws.Range("I134").Value = 1392.4286
$ws.Range("J134").Value = 3156.75
$ws.Range("K134").Value = 4177.2858
$ws.Range("L134").Value = 9470.25
$ws.Range("M134").Value = -1642.2858
$ws.Range("N134").Value = -14540.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2263.8823
$ws.Range("I31").Value = 1478.88
$ws.Range("J31").Value = 4444.4443
$ws.Range("K31").Value = 1478.88
$ws.Range("L31").Value = 4444.4443
$ws.Range("M31").Value = -1183.88
$ws.Range("N31").Value = -5034.4443
$ws.Range("H34").Value = 2263.8823
$ws.Range("I34").Value = 1478.88
$ws.Range("J34").Value = 4444.4443
$ws.Range("K34").Value = 1478.88
$ws.Range("L34").Value = 4444.4443
$ws.Range("M34").Value = -1276.88
$ws.Range("N34").Value = -4848.4443
$ws.Range("H58").Value = 904682.5600000001
$ws.Range("I58").Value = 1425687.1
$ws.Range("K58").Value = 1425687.1
$ws.Range("M58").Value = -1425484.1
$ws.Range("H132").Value = 234100.69
$ws.Range("I132").Value = 308007.53
$ws.Range("J132").Value = 1822.0714
$ws.Range("K132").Value = 924022.5900000001
$ws.Range("L132").Value = 5466.2142
$ws.Range("M132").Value = -921492.5900000001
$ws.Range("N132").Value = -10526.2142
$ws.Range("H134").Value = 1173.4722
$ws.Range("I134").Value = 1006.0645
$ws.Range("J134").Value = 2211.4
$ws.Range("K134").Value = 3018.1935
$ws.Range("L134").Value = 6634.200000000001
$ws.Range("M134").Value = -483.1934999999999
$ws.Range("N134").Value = -11704.2
$ws.Range("H136").Value = 904682.5600000001
$ws.Range("I136").Value = 1425687.1
$ws.Range("K136").Value = 4277061.300000001
$ws.Range("M136").Value = -4274511.300000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 7703.2
$ws.Range("J75").Value = 11750
$ws.Range("L75").Value = 35250
$ws.Range("N75").Value = -37246
$ws.Range("H78").Value = 7703.2
$ws.Range("J78").Value = 11750
$ws.Range("L78").Value = 105750
$ws.Range("N78").Value = -115734
$ws.Range("H107").Value = 446.35715
$ws.Range("I107").Value = 411.46155
$ws.Range("J107").Value = 900
$ws.Range("K107").Value = 1234.38465
$ws.Range("L107").Value = 2700
$ws.Range("M107").Value = 685.61535
$ws.Range("N107").Value = -6540
$ws.Range("H122").Value = 861.619
$ws.Range("I122").Value = 665.75
$ws.Range("J122").Value = 1122.7778
$ws.Range("K122").Value = 5991.75
$ws.Range("L122").Value = 10105.0002
$ws.Range("M122").Value = -3541.75
$ws.Range("N122").Value = -15005.0002
$ws.Range("H131").Value = 17260320
$ws.Range("I131").Value = 15928.571
$ws.Range("J131").Value = 19627196
$ws.Range("K131").Value = 47785.713
$ws.Range("L131").Value = 58881588
$ws.Range("M131").Value = -42745.713
$ws.Range("N131").Value = -58891668

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1144.9166
$ws.Range("I132").Value = 719.9655
$ws.Range("K132").Value = 2159.8965
$ws.Range("M132").Value = 370.1035000000002
$ws.Range("H138").Value = 44429
$ws.Range("J138").Value = 44429
$ws.Range("L138").Value = 44429
$ws.Range("N138").Value = -54709

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3350.745
$ws.Range("I132").Value = 2908.2104
$ws.Range("J132").Value = 4644.3076
$ws.Range("K132").Value = 8724.6312
$ws.Range("L132").Value = 13932.9228
$ws.Range("M132").Value = -6194.6312
$ws.Range("N132").Value = -18992.9228

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4199.6
$ws.Range("I62").Value = 3999.5
$ws.Range("J62").Value = 5000
$ws.Range("K62").Value = 3999.5
$ws.Range("L62").Value = 5000
$ws.Range("M62").Value = -3375.5
$ws.Range("N62").Value = -6248
$ws.Range("H65").Value = 4199.6
$ws.Range("I65").Value = 3999.5
$ws.Range("J65").Value = 5000
$ws.Range("K65").Value = 19997.5
$ws.Range("L65").Value = 25000
$ws.Range("M65").Value = -16877.5
$ws.Range("N65").Value = -31240
$ws.Range("H80").Value = 28767
$ws.Range("J80").Value = 28767
$ws.Range("L80").Value = 28767
$ws.Range("N80").Value = -30763
$ws.Range("H83").Value = 28767
$ws.Range("J83").Value = 28767
$ws.Range("L83").Value = 86301
$ws.Range("N83").Value = -96285
$ws.Range("H122").Value = 8447677
$ws.Range("I122").Value = 10418169
$ws.Range("J122").Value = 4809846
$ws.Range("K122").Value = 31254507
$ws.Range("L122").Value = 14429538
$ws.Range("M122").Value = -31252057
$ws.Range("N122").Value = -14434438
$ws.Range("H132").Value = 1446.9412
$ws.Range("I132").Value = 967.4231
$ws.Range("J132").Value = 3005.375
$ws.Range("K132").Value = 2902.2693
$ws.Range("L132").Value = 9016.125
$ws.Range("M132").Value = -372.2692999999999
$ws.Range("N132").Value = -14076.125
$ws.Range("H136").Value = 1588.9803
$ws.Range("I136").Value = 1309.45
$ws.Range("J136").Value = 2605.4546
$ws.Range("K136").Value = 3928.35
$ws.Range("L136").Value = 7816.3638
$ws.Range("M136").Value = -1378.35
$ws.Range("N136").Value = -12916.3638
$ws.Range("H138").Value = 82585.8
$ws.Range("J138").Value = 82585.8
$ws.Range("L138").Value = 82585.8
$ws.Range("N138").Value = -92865.8
